# Commit message: "Added reels into chgSymbols"
# The data rows 2-21 (columns A:F) of the active sheet are reshuffled
# to a new order (same set of rows, different row positions).
# We rewrite the values of A2:F21 directly to match the target layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(701,3,90,45,97,15),
    @(501,9,52,30,75,45),
    @(1001,18,30,75,60,72),
    @(201,9,30,15,45,30),
    @(902,1,0,0,0,0),
    @(401,9,48,67,75,45),
    @(601,9,60,67,60,42),
    @(1201,2,10,10,10,10),
    @(1203,3,15,15,15,15),
    @(901,16,15,45,60,60),
    @(301,6,45,30,60,45),
    @(801,3,67,65,52,45),
    @(1202,2,10,10,10,10),
    @(101,9,30,15,60,15),
    @(1101,0,15,30,30,0),
    @(2,0,2,2,2,2),
    @(3,0,3,3,3,3),
    @(1,0,2,2,2,2),
    @(502,0,4,0,0,0),
    @(802,0,4,5,4,0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $rowValues[$c - 1]
    }
}
